# Cuttings log well 16-2-7 — classify "Marl" and "Tuff" lithologies as
# "Shale" in the simplified-lithology lookup formula (column D), and move
# the active selection the way the author left it (scrolled down, cell E6
# selected) after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D2: first (non-shared) occurrence of the lookup formula ---------------
$ws.Range("D2").Formula = '=IF(OR(C2="Claystone",C2="Siltstone",C2="Shale",C2="Clay",C2="Silty clay-shale",C2="Clay-shale",C2="Marl",C2="Tuff"),"Shale",IF(C2="Limestone","Carbonate",C2))'

# --- D3:D69: shared formula block, same edit (adds Marl/Tuff to the OR) ----
$ws.Range("D3:D69").Formula = '=IF(OR(C3="Claystone",C3="Siltstone",C3="Shale",C3="Clay",C3="Silty clay-shale",C3="Clay-shale",C3="Marl",C3="Tuff"),"Shale",IF(C3="Limestone","Carbonate",C3))'

# --- View state: scroll down so row 58 is at the top, then select E6 -------
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("E6").Select() | Out-Null
